$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes look like plain numbers (e.g. "494.11") even though
# the sheet stores every Price/Volume value as text (note multi-dot values like
# "54.447.26" which are not valid numbers). Force text interpretation via the
# NumberFormat trick, then restore the default "Normal" style so the cell style
# index matches the original (unstyled) cells.
function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '54.447.26'
$ws.Range('E2').Value = '  -2.51%  '
Set-TextValue 'D3' '2.288.46'
$ws.Range('E3').Value = '  -2.86%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '494.11'
$ws.Range('E5').Value = '  -1.98%  '
Set-TextValue 'D6' '127.28'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  -1.63%  '
Set-TextValue 'D9' '2.288.58'
$ws.Range('E9').Value = '  -3.42%  '
Set-TextValue 'D10' '0.0948'
$ws.Range('E10').Value = '  -2.36%  '
Set-TextValue 'D11' '0.150'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E13').Value = '  -3.19%  '
Set-TextValue 'D14' '2.694.56'
Set-TextValue 'D15' '21.57'
$ws.Range('E15').Value = '  +0.52%  '
Set-TextValue 'D16' '54.375.75'
$ws.Range('E16').Value = '  -2.53%  '
$ws.Range('E17').Value = '  -2.04%  '
Set-TextValue 'D18' '2.275.44'
$ws.Range('E18').Value = '  -4.06%  '
Set-TextValue 'D19' '10.00'
$ws.Range('E19').Value = '  +0.99%  '
$ws.Range('E20').Value = '  +1.03%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D21' '303.75'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '6.48'
$ws.Range('E22').Value = '  +4.35%  '
$ws.Range('E23').Value = '  +0.42%  '
Set-TextValue 'D24' '5.39'
$ws.Range('E24').Value = '  -2.80%  '
Set-TextValue 'D25' '63.50'
$ws.Range('E25').Value = '  -2.71%  '
Set-TextValue 'D26' '0.999'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  +1.08%  '
Set-TextValue 'D28' '2.395.09'
$ws.Range('E28').Value = '  -2.91%  '
$ws.Range('E29').Value = '  +2.77%  '
$ws.Range('E30').Value = '  -0.16%  '
Set-TextValue 'D31' '169.25'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('E32').Value = '  -2.15%  '
Set-TextValue 'D33' '0.0₃0685'
$ws.Range('E33').Value = '  -3.02%  '
Set-TextValue 'D34' '5.88'
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +1.63%  '
Set-TextValue 'D38' '17.62'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('E39').Value = '  +2.60%  '
Set-TextValue 'D40' '0.870'
$ws.Range('E40').Value = '  +3.97%  '
$ws.Range('E41').Value = '  +0.03%  '
Set-TextValue 'D42' '35.52'
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('E43').Value = '  +1.31%  '
$ws.Range('E44').Value = '  +0.54%  '
Set-TextValue 'D45' '3.35'
$ws.Range('E45').Value = '  +0.22%  '
Set-TextValue 'D46' '128.64'
$ws.Range('E46').Value = '  +2.25%  '
Set-TextValue 'D47' '4.80'
$ws.Range('E47').Value = '  -1.07%  '
Set-TextValue 'D49' '0.544'
$ws.Range('E49').Value = '  -2.39%  '
Set-TextValue 'D50' '239.32'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('E51').Value = '  +0.30%  '
